# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 69 (pushing the existing rows 69-80
# down to 70-81) and populate it with the new week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 69:80 down to 70:81, leaving row 69 free for the new data.
$ws.Rows.Item(69).Insert()

$ws.Range("A69").Value = 10
$ws.Range("B69").Value = "Vega Modelo de Temuco"
$ws.Range("C69").Value = "La Araucanía"
$ws.Range("D69").Value = 44637
$ws.Range("E69").Value = 9
$ws.Range("F69").Value = 100112030
$ws.Range("G69").Value = "Poroto granado"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 50
$ws.Range("K69").Value = 20000
$ws.Range("L69").Value = 20000
$ws.Range("M69").Value = 20000
$ws.Range("N69").Value = "$/saco 25 kilos"
$ws.Range("O69").Value = "Región de La Araucanía"
$ws.Range("P69").Value = 800
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"

Write-Output "Inserted new row 69; rows 69-80 shifted to 70-81."
